# Swap the order of "System" and "dnasr281@gmail.com" in the
# "Recorded By" column (G) wherever the cell contains exactly those
# two entries, toggling between "System, dnasr281@gmail.com" and
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nameA = "System, dnasr281@gmail.com"
$nameB = "dnasr281@gmail.com, System"

$lastRow = $ws.Cells($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $nameA) {
        $cell.Value = $nameB
    } elseif ($val -eq $nameB) {
        $cell.Value = $nameA
    }
}
